$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newIds = @(
    "PEP_ID-2009336",
    "PEP_ID-2009337",
    "PEP_ID-2009338",
    "PEP_ID-2009339",
    "PEP_ID-2009340",
    "PEP_ID-2009342",
    "PEP_ID-2009343",
    "PEP_ID-2009344",
    "PEP_ID-2009347",
    "PEP_ID-2009349",
    "PEP_ID-2009350",
    "PEP_ID-2009351",
    "PEP_ID-2009352",
    "PEP_ID-2009353",
    "PEP_ID-2009354",
    "PEP_ID-2009357",
    "PEP_ID-2009358"
)

$startRow = 12
for ($i = 0; $i -lt $newIds.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newIds[$i]
}
